# Updates cryptos list prices / 1h volume percentages (and reorders a
# couple of rows) to match the latest scrape, per the GitHub Actions
# commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that looks like a plain number (e.g. "1.00",
# "62.15") while keeping the cell as TEXT, matching the source data
# (these are inline/shared strings, not numeric cells). We briefly force
# a text number format, set the value, then restore the default "Normal"
# style so no stray formatting differences are introduced.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.978.26"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.561.36"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("E9").Value = "  +0.29%  "
Set-TextValue "D10" "0.0598"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.783.16"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.562.61"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +0.80%  "
Set-TextValue "D16" "62.15"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "26.990.36"
$ws.Range("E17").Value = "  +0.27%  "
Set-TextValue "D18" "217.26"
$ws.Range("E18").Value = "  +0.05%  "
Set-TextValue "D20" "7.37"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -1.24%  "
Set-TextValue "D25" "153.42"
$ws.Range("E25").Value = "  -0.29%  "
Set-TextValue "D26" "6.62"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  +1.46%  "
Set-TextValue "D29" "1.01"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +0.74%  "
Set-TextValue "D31" "1.12"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.424.68"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "3.11"
$ws.Range("E34").Value = "  +3.61%  "
Set-TextValue "D35" "1.60"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  +8.90%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("E38").Value = "  +0.73%  "
Set-TextValue "D39" "0.535"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("E44").Value = "  +1.78%  "
Set-TextValue "D45" "64.83"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "1.697.04"
$ws.Range("E47").Value = "  +0.50%  "
Set-TextValue "D48" "87.42"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0957"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  -0.23%  "
